$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 15, shifting existing rows 15-78 down to 16-79.
$ws.Rows("15:15").Insert(-4121)

# Populate the newly inserted row 15 with the new weekly price entry.
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C15").Value = "Arica y Parinacota"
$ws.Range("D15").Value = 44707
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = 100112021
$ws.Range("G15").Value = "Ají"
$ws.Range("H15").Value = "Inferno"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 120
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14500
$ws.Range("N15").Value = "$/caja 15 kilos"
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 967
$ws.Range("Q15").Value = 15
$ws.Range("R15").Value = "Hortaliza"
